$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each run of the watched process dropped a timestamped "first chars" entry:
# [timestamp, first-char-seen, elapsed-seconds]. Write them in, one row per run.
$entries = @(
    @("Mon_Dec__4_21_42_49_2023", "f", 30),
    @("Mon_Dec__4_21_47_16_2023", "f", 30),
    @("Mon_Dec__4_22_14_38_2023", "f", 30)
)

$row = 1
foreach ($entry in $entries) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}

# Animate/highlight the freshly written block for review.
$ws.Range("A1:C9").Select() | Out-Null
